$wb = $excel.ActiveWorkbook

# --- "param" sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("param")

# c_D / u_D distribution: "NA" -> "fixed"
$ws.Range("C10").Value = "fixed"
$ws.Range("C16").Value = "fixed"

# par1/par2/par3 numeric columns: custom "0.0000" number format -> built-in "0.00"
$ws.Range("D2:F6").NumberFormat = "0.00"
$ws.Range("D7:F7").NumberFormat = "0.00"
$ws.Range("D8:F12").NumberFormat = "0.00"
$ws.Range("D13:F13").NumberFormat = "0.00"
$ws.Range("D14:F16").NumberFormat = "0.00"

# the "distribution" dropdown validation list shrank by one row (B2:B10 -> B2:B9)
$ws.Range("C2:C17").Validation.Delete()
$ws.Range("C2:C17").Validation.Add(3, 1, 1, "='distribution names'!`$B`$2:`$B`$9")

$ws.Range("C16").Select()

# --- "distribution names" sheet --------------------------------------------
$ws2 = $wb.Worksheets.Item("distribution names")

# Type column shifts up: weibull row becomes "uniform", uniform row becomes
# "fixed", and the trailing "dirichlet " / "NA" rows are dropped entirely.
$ws2.Range("B7").Value = "uniform"
$ws2.Range("B8").Value = "fixed"
$ws2.Range("B9").Value = ""
$ws2.Range("B10").Value = ""

$ws2.Range("B9").Select()
